$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): B11 4 -> 5 ; C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): B12 84 -> 105 ; C12 -3 -> -3.6 ; E12 "81/112" -> "101.4/140"
$ws.Range("B12").Value = 105
$ws.Range("C12").Value = -3.6
$ws.Range("E12").Value = "101.4/140"
